# GEORGIA_2020.xlsx cleanup script
# - rename header columns to snake_case
# - title-case Spanish connector words (de/del/el/la/las/los/y) that were
#   left lowercase in municipality / state names
# - fix a couple of one-off casing glitches (GUANAJUATO, MonteMorelos)
# - drop the trailing footnote rows (no longer part of the clean table)
# - nudge a handful of percentage cells to the freshly recomputed values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function TitleConnectors([string]$s) {
    $connectors = @("de", "del", "el", "la", "las", "los", "y")
    $words = $s.Split(" ")
    $out = @()
    foreach ($w in $words) {
        $isConnector = $false
        foreach ($c in $connectors) {
            if ($w.Equals($c)) {
                $isConnector = $true
            }
        }
        if ($isConnector) {
            $out += ($w.Substring(0, 1).ToUpper() + $w.Substring(1))
        } else {
            $out += $w
        }
    }
    return ($out -join " ")
}

# 1) Drop the footnote rows (1385:1390) - they are free text, not data.
$ws.Rows("1385:1390").Delete()

# 2) Rename the header row to the new snake_case column names.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 3) Title-case the lowercase Spanish connector words throughout columns A/B.
$lastRow = 1384
for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in @("A", "B")) {
        $cell = $ws.Range($col + $r)
        $v = $cell.Value2
        if ($v -ne $null -and $v -is [string]) {
            $newV = TitleConnectors $v
            if (-not $newV.Equals($v)) {
                $cell.Value = $newV
            }
        }
    }
}

# 4) One-off casing fixes that aren't covered by the connector-word rule.
$ws.Range("A291").Value = "Guanajuato"
$ws.Range("B708").Value = "Montemorelos"

# 5) A handful of percentage cells were recomputed and differ in the last
#    significant digit - set them to the refreshed values explicitly.
$rowsWithC14 = @(30, 163, 257, 422, 461, 536, 639, 650, 715, 1071, 1090, 1173, 1191, 1223, 1258, 1286, 1352, 1382)
foreach ($r in $rowsWithC14) {
    $ws.Range("D" + $r).Value = 0.0009155113784985612
}

$rowsWithC15 = @(371, 379, 697, 828, 1270)
foreach ($r in $rowsWithC15) {
    $ws.Range("D" + $r).Value = 0.0009809050483913155
}

$ws.Range("D334").Value = 0.09135495684017789

# 6) Make sure the declared dimension matches the trimmed data range.
$ws.Range("A1").Worksheet.UsedRange | Out-Null
